$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.977.04'
$ws.Range("E2").Value = '  +1.82%  '
$ws.Range("D3").Value = '1.673.40'
$ws.Range("E3").Value = '  +0.93%  '
$ws.Range("E4").Value = '  -0.22%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '331.10'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +7.75%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.9990'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.18%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.3655'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +1.25%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '47.26'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -0.48%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.3226'
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '1.143'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +2.32%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.07129'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +2.28%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.9992'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -0.13%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '6.080'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +3.61%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '19.64'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("D15").Value = '1.667.01'
$ws.Range("E15").Value = '  +1.17%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '6.656'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +1.78%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.00001046'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +0.41%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '0.06547'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +0.38%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.9992'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -0.05%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '78.88'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +3.45%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '15.82'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +1.30%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '5.909'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.20%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '12.79'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("D24").Value = '24.970.73'
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.437'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -0.99%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.386'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +4.54%  '
$ws.Range("E27").Value = '  +1.49%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '18.68'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +1.55%  '
$ws.Range("D29").Value = '1.851.46'
$ws.Range("E29").Value = '  +1.03%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '125.85'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +1.88%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '1.195'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +0.57%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '4.079'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -0.03%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '5.787'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +2.31%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.08450'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +1.23%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.671'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -1.66%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '12.30'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -0.19%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '5.159'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.54%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.06052'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.15%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '1.229'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +2.21%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.02229'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +1.79%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.2089'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +2.00%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '8.235'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +0.37%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.9985'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.22%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.5956'
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '13.64'
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '3.838'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +2.74%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.5725'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +2.74%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '123.94'
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '1.960'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +1.32%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.07005'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +1.31%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.193'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +3.67%  '
